# Commit: Re-do TODE grade norms with week coding for grade
#
# The raw-score -> scale-score ("ss") lookup tables on each of the six
# grade/season tabs (K-Fall, K-Spring, 1-Fall, 1-Spring, 2-Fall, 2-Spring)
# were recomputed. This script rewrites column B (the "ss" values) on
# each sheet to the newly recomputed lookup values. Column A (the raw
# score / row index) and the header row are left untouched.

$wb = $excel.ActiveWorkbook

# --- Sheet: K-Fall (19 updated cell(s) in column B) ---
$ws = $wb.Worksheets.Item("K-Fall")
$rows   = @(3,4,9,15,38,39,44,45,50,51,56,57,62,63,68,69,74,75,76)
$values = @(66,67,71,76,96,97,101,102,106,107,111,112,116,117,121,122,126,127,128)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 2).Value = $values[$i]
}

# --- Sheet: K-Spring (7 updated cell(s) in column B) ---
$ws = $wb.Worksheets.Item("K-Spring")
$rows   = @(3,4,8,12,17,21,30)
$values = @(61,62,65,68,72,75,82)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 2).Value = $values[$i]
}

# --- Sheet: 1-Fall (20 updated cell(s) in column B) ---
$ws = $wb.Worksheets.Item("1-Fall")
$rows   = @(2,3,6,10,13,17,21,24,28,32,72,77,81,86,87,91,92,96,97,98)
$values = @(56,57,59,62,64,67,70,72,75,78,108,112,115,119,120,123,124,127,128,129)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 2).Value = $values[$i]
}

# --- Sheet: 1-Spring (103 updated cell(s) in column B) ---
$ws = $wb.Worksheets.Item("1-Spring")
$rows   = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,71,72,73,74,75,76,77,78,79,80,81,82,83,84,85,86,87,88,89,90,91,92,93,94,95,96,97,98,99,100,101,102,103,104)
$values = @(66,66,67,68,69,70,71,71,72,73,74,75,76,76,77,78,79,80,81,81,82,83,84,85,86,86,87,88,89,90,91,91,92,93,94,95,96,97,97,98,99,100,101,102,102,103,104,105,106,107,107,108,109,110,111,112,112,113,114,115,116,117,117,118,119,120,121,122,122,123,124,125,126,127,128,128,129,130,130,130,130,130,130,130,130,130,130,130,130,130,130,130,130,130,130,130,130,130,130,130,130,130,130)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 2).Value = $values[$i]
}

# --- Sheet: 2-Fall (11 updated cell(s) in column B) ---
$ws = $wb.Worksheets.Item("2-Fall")
$rows   = @(3,6,8,11,14,17,20,94,102,103,104)
$values = @(50,52,53,55,57,59,61,114,121,122,123)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 2).Value = $values[$i]
}

# --- Sheet: 2-Spring (20 updated cell(s) in column B) ---
$ws = $wb.Worksheets.Item("2-Spring")
$rows   = @(3,6,8,11,13,16,21,24,29,40,61,64,71,75,79,83,88,94,95,112)
$values = @(47,49,50,52,53,55,58,60,63,70,85,87,92,95,98,101,105,110,111,128)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 2).Value = $values[$i]
}
